# Apply cryptocurrency price/volume updates per commit:
# "Updated cryptos list on Mon Aug 14 23:39:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    # Force the value to be stored as literal text (not auto-converted to a
    # number/date by Excel), while leaving the cell's style index untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.390.57'
Set-TextValue $ws.Range('E2') '  +0.39%  '
Set-TextValue $ws.Range('D3') '1.843.26'
Set-TextValue $ws.Range('E3') '  +0.21%  '
Set-TextValue $ws.Range('D4') '0.9990'
Set-TextValue $ws.Range('E4') '  -0.05%  '
Set-TextValue $ws.Range('D5') '240.24'
Set-TextValue $ws.Range('E5') '  +0.17%  '
Set-TextValue $ws.Range('D6') '0.6342'
Set-TextValue $ws.Range('E6') '  +1.30%  '
Set-TextValue $ws.Range('D7') '0.9997'
Set-TextValue $ws.Range('E7') '  +0.12%  '
Set-TextValue $ws.Range('E8') '  +0.02%  '
Set-TextValue $ws.Range('D9') '25.12'
Set-TextValue $ws.Range('E9') '  +3.46%  '
Set-TextValue $ws.Range('D10') '0.2905'
Set-TextValue $ws.Range('E10') '  +0.50%  '
Set-TextValue $ws.Range('D11') '0.07745'
Set-TextValue $ws.Range('E11') '  +0.53%  '
Set-TextValue $ws.Range('D12') '1.902.59'
Set-TextValue $ws.Range('E12') '  +3.44%  '
Set-TextValue $ws.Range('D13') '4.990'
Set-TextValue $ws.Range('E13') '  +0.16%  '
Set-TextValue $ws.Range('D14') '0.6800'
Set-TextValue $ws.Range('E14') '  +0.49%  '
Set-TextValue $ws.Range('D15') '0.00001026'
Set-TextValue $ws.Range('E15') '  -0.09%  '
Set-TextValue $ws.Range('D16') '82.02'
Set-TextValue $ws.Range('E16') '  +0.02%  '
Set-TextValue $ws.Range('D17') '6.267'
Set-TextValue $ws.Range('E17') '  +2.86%  '
Set-TextValue $ws.Range('D18') '29.389.01'
Set-TextValue $ws.Range('E18') '  +0.34%  '
Set-TextValue $ws.Range('D19') '230.51'
Set-TextValue $ws.Range('D20') '12.34'
Set-TextValue $ws.Range('E20') '  +0.87%  '
Set-TextValue $ws.Range('E21') '  +0.04%  '
Set-TextValue $ws.Range('D22') '7.427'
Set-TextValue $ws.Range('E22') '  +0.80%  '
Set-TextValue $ws.Range('D23') '1.000'
Set-TextValue $ws.Range('E23') '  +0.16%  '
Set-TextValue $ws.Range('D24') '158.05'
Set-TextValue $ws.Range('E24') '  -0.21%  '
Set-TextValue $ws.Range('D25') '8.508'
Set-TextValue $ws.Range('E25') '  +1.74%  '
Set-TextValue $ws.Range('D26') '0.1359'
Set-TextValue $ws.Range('E26') '  -1.52%  '
Set-TextValue $ws.Range('D27') '17.49'
Set-TextValue $ws.Range('E27') '  -0.17%  '
Set-TextValue $ws.Range('D28') '0.06572'
Set-TextValue $ws.Range('E28') '  +15.35%  '
Set-TextValue $ws.Range('D29') '1.431'
Set-TextValue $ws.Range('E29') '  +2.70%  '
Set-TextValue $ws.Range('D30') '1.490'
Set-TextValue $ws.Range('E30') '  +1.24%  '
Set-TextValue $ws.Range('D31') '4.077'
Set-TextValue $ws.Range('E31') '  -0.40%  '
Set-TextValue $ws.Range('D32') '4.053'
Set-TextValue $ws.Range('E32') '  +0.65%  '
Set-TextValue $ws.Range('D33') '1.839'
Set-TextValue $ws.Range('E33') '  +1.26%  '
Set-TextValue $ws.Range('D34') '1.142'
Set-TextValue $ws.Range('E34') '  +0.13%  '
Set-TextValue $ws.Range('D35') '0.7011'
Set-TextValue $ws.Range('E35') '  +0.08%  '
Set-TextValue $ws.Range('D36') '2.578'
Set-TextValue $ws.Range('E36') '  -0.09%  '
Set-TextValue $ws.Range('D37') '0.01859'
Set-TextValue $ws.Range('E37') '  +2.68%  '
Set-TextValue $ws.Range('D38') '1.250.47'
Set-TextValue $ws.Range('E38') '  +0.65%  '
Set-TextValue $ws.Range('D39') '2.819'
Set-TextValue $ws.Range('D40') '6.770'
Set-TextValue $ws.Range('E40') '  +4.04%  '
Set-TextValue $ws.Range('D41') '0.9381'
Set-TextValue $ws.Range('E41') '  +3.90%  '
Set-TextValue $ws.Range('D42') '0.9995'
Set-TextValue $ws.Range('E42') '  +0.23%  '
Set-TextValue $ws.Range('D43') '2.003.15'
Set-TextValue $ws.Range('E43') '  +0.12%  '
Set-TextValue $ws.Range('D44') '101.16'
Set-TextValue $ws.Range('E44') '  -0.01%  '
Set-TextValue $ws.Range('D45') '65.46'
Set-TextValue $ws.Range('E45') '  -0.30%  '
Set-TextValue $ws.Range('B46') 'BabyDogeCoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.00000000119'
Set-TextValue $ws.Range('E46') '  +5.69%  '
Set-TextValue $ws.Range('B47') 'Aptos'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D47') '7.071'
Set-TextValue $ws.Range('E47') '  +0.08%  '
Set-TextValue $ws.Range('B48') 'RenderToken'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D48') '1.721'
Set-TextValue $ws.Range('E48') '  +4.11%  '
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '9.057'
Set-TextValue $ws.Range('E49') '  +0.99%  '
Set-TextValue $ws.Range('B50') 'Algorand'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D50') '0.1150'
Set-TextValue $ws.Range('E50') '  -1.41%  '
Set-TextValue $ws.Range('B51') 'TheSandbox'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D51') '0.3917'
Set-TextValue $ws.Range('E51') '  -0.34%  '
